$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

$ws1.Range("B1").Value = "2485-RBI-EPP-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-2-DATE-VAR-INST-NO-1st"
$ws1.Range("B2").Value = "248e"
$ws2.Range("B1").Value = "2485-RBI-EPP-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-2-DATE-VAR-INST-NO-1st"

$ws1.Range("B1").Select()
$ws2.Select()
